# Auto-generated edit script applying numeric corrections to Cactuar_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 12666.637
$ws.Range("J40").Value = 11515.429
$ws.Range("L40").Value = 11515.429
$ws.Range("N40").Value = -11865.429

$ws.Range("H43").Value = 3075.8333
$ws.Range("J43").Value = 3291
$ws.Range("L43").Value = 3291
$ws.Range("N43").Value = -3429

$ws.Range("H80").Value = 43104948
$ws.Range("J80").Value = 66667860
$ws.Range("L80").Value = 200003580
$ws.Range("N80").Value = -200005576

$ws.Range("H83").Value = 43104948
$ws.Range("J83").Value = 66667860
$ws.Range("L83").Value = 600010740
$ws.Range("N83").Value = -600020724

$ws.Range("H88").Value = 9308.352999999999
$ws.Range("J88").Value = 12832.75
$ws.Range("L88").Value = 12832.75
$ws.Range("N88").Value = -13644.75

$ws.Range("H91").Value = 9308.352999999999
$ws.Range("J91").Value = 12832.75
$ws.Range("L91").Value = 12832.75
$ws.Range("N91").Value = -15640.75

$ws.Range("H103").Value = 604.375
$ws.Range("I103").Value = 664
$ws.Range("J103").Value = 568.6
$ws.Range("K103").Value = 1992
$ws.Range("L103").Value = 1705.8
$ws.Range("M103").Value = -1406
$ws.Range("N103").Value = -2877.8

$ws.Range("H133").Value = 94682.47
$ws.Range("J133").Value = 94682.47
$ws.Range("L133").Value = 94682.47
$ws.Range("N133").Value = -104802.47

$ws.Range("H137").Value = 8810779
$ws.Range("I137").Value = 460757.4
$ws.Range("J137").Value = 19616690
$ws.Range("K137").Value = 1382272.2
$ws.Range("L137").Value = 58850070
$ws.Range("M137").Value = -1379722.2
$ws.Range("N137").Value = -58855170

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 43332.668
$ws.Range("J76").Value = 43332.668
$ws.Range("L76").Value = 43332.668
$ws.Range("N76").Value = -44008.668

$ws.Range("H79").Value = 43332.668
$ws.Range("J79").Value = 43332.668
$ws.Range("L79").Value = 43332.668
$ws.Range("N79").Value = -45672.668

$ws.Range("H110").Value = 1313.174
$ws.Range("I110").Value = 1133.1765
$ws.Range("J110").Value = 1823.1666
$ws.Range("K110").Value = 1133.1765
$ws.Range("L110").Value = 1823.1666
$ws.Range("M110").Value = 911.8235
$ws.Range("N110").Value = -5913.1666

$ws.Range("H122").Value = 3259.3333
$ws.Range("I122").Value = 1745.1
$ws.Range("J122").Value = 7585.7144
$ws.Range("K122").Value = 5235.299999999999
$ws.Range("L122").Value = 22757.1432
$ws.Range("M122").Value = -2785.299999999999
$ws.Range("N122").Value = -27657.1432

$ws.Range("H132").Value = 4293.4287
$ws.Range("I132").Value = 1573.2059
$ws.Range("J132").Value = 10459.267
$ws.Range("K132").Value = 4719.6177
$ws.Range("L132").Value = 31377.801
$ws.Range("M132").Value = -2189.6177
$ws.Range("N132").Value = -36437.801

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1840.3889
$ws.Range("I94").Value = 1897
$ws.Range("K94").Value = 1897
$ws.Range("M94").Value = -1446

$ws.Range("H107").Value = 3895.2
$ws.Range("I107").Value = 2619.6875
$ws.Range("K107").Value = 2619.6875
$ws.Range("M107").Value = -699.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2963.375
$ws.Range("I16").Value = 1956.2858
$ws.Range("K16").Value = 1956.2858
$ws.Range("M16").Value = -1669.2858

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H35").Value = 1559.8
$ws.Range("I35").Value = 1599.6666
$ws.Range("K35").Value = 1599.6666
$ws.Range("M35").Value = -1305.6666

$ws.Range("H94").Value = 2198.5454
$ws.Range("I94").Value = 305.25
$ws.Range("J94").Value = 2619.2778
$ws.Range("K94").Value = 305.25
$ws.Range("L94").Value = 2619.2778
$ws.Range("M94").Value = 145.75
$ws.Range("N94").Value = -3521.2778

$ws.Range("H113").Value = 2963.375
$ws.Range("I113").Value = 1956.2858
$ws.Range("K113").Value = 1956.2858
$ws.Range("M113").Value = 213.7141999999999

$ws.Range("H122").Value = 4010.7646
$ws.Range("I122").Value = 2561.1538
$ws.Range("K122").Value = 7683.4614
$ws.Range("M122").Value = -5233.4614

$ws.Range("H132").Value = 44446084
$ws.Range("I132").Value = 51283212
$ws.Range("K132").Value = 153849636
$ws.Range("M132").Value = -153847106

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 671.3333
$ws.Range("J113").Value = 611.1818
$ws.Range("L113").Value = 1833.5454
$ws.Range("N113").Value = -6173.5454

$ws.Range("H137").Value = 53573996
$ws.Range("I137").Value = 93751490
$ws.Range("J137").Value = 4004.5
$ws.Range("K137").Value = 281254470
$ws.Range("L137").Value = 12013.5
$ws.Range("M137").Value = -281249370
$ws.Range("N137").Value = -22213.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6493.3076
$ws.Range("I70").Value = 6129
$ws.Range("K70").Value = 6129
$ws.Range("M70").Value = -5859

$ws.Range("H73").Value = 6493.3076
$ws.Range("I73").Value = 6129
$ws.Range("K73").Value = 6129
$ws.Range("M73").Value = -5193

$ws.Range("H113").Value = 2314.9167
$ws.Range("I113").Value = 1417.1111
$ws.Range("K113").Value = 1417.1111
$ws.Range("M113").Value = 752.8888999999999

$ws.Range("H122").Value = 594948.5600000001
$ws.Range("I122").Value = 5002499
$ws.Range("J122").Value = 7275.2
$ws.Range("K122").Value = 15007497
$ws.Range("L122").Value = 21825.6
$ws.Range("M122").Value = -15005047
$ws.Range("N122").Value = -26725.6

$ws.Range("H132").Value = 77457.55499999999
$ws.Range("I132").Value = 98286.05
$ws.Range("J132").Value = 4557.8335
$ws.Range("K132").Value = 294858.15
$ws.Range("L132").Value = 13673.5005
$ws.Range("M132").Value = -292328.15
$ws.Range("N132").Value = -18733.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8750
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 8750
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 8750
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -8974

$ws.Range("H22").Value = 975
$ws.Range("I22").Value = 849.8182
$ws.Range("K22").Value = 849.8182
$ws.Range("M22").Value = -554.8182

$ws.Range("H27").Value = 975
$ws.Range("I27").Value = 849.8182
$ws.Range("K27").Value = 849.8182
$ws.Range("M27").Value = -742.8182

$ws.Range("H46").Value = 4292.7915
$ws.Range("I46").Value = 3996.6667
$ws.Range("K46").Value = 3996.6667
$ws.Range("M46").Value = -3808.6667

$ws.Range("H50").Value = 42584
$ws.Range("J50").Value = 42584
$ws.Range("L50").Value = 42584
$ws.Range("N50").Value = -43858

$ws.Range("H55").Value = 553.7692
$ws.Range("J55").Value = 927.7143
$ws.Range("L55").Value = 927.7143
$ws.Range("N55").Value = -1273.7143

$ws.Range("H126").Value = 8750
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 8750
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 26250
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -31190

$ws.Range("H132").Value = 3716.71
$ws.Range("I132").Value = 3208.7693
$ws.Range("J132").Value = 4660.029
$ws.Range("K132").Value = 9626.3079
$ws.Range("L132").Value = 13980.087
$ws.Range("M132").Value = -7096.3079
$ws.Range("N132").Value = -19040.087

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H52").Value = 7750
$ws.Range("I52").Value = 7750
$ws.Range("K52").Value = 7750
$ws.Range("M52").Value = -7524

$ws.Range("H54").Value = 13043.889
$ws.Range("I54").Value = 5774.2856
$ws.Range("J54").Value = 38487.5
$ws.Range("K54").Value = 5774.2856
$ws.Range("L54").Value = 38487.5
$ws.Range("M54").Value = -5254.2856
$ws.Range("N54").Value = -39527.5

$ws.Range("H81").Value = 6568
$ws.Range("I81").Value = 5296.5
$ws.Range("J81").Value = 7839.5
$ws.Range("K81").Value = 10593
$ws.Range("L81").Value = 15679
$ws.Range("M81").Value = -9532
$ws.Range("N81").Value = -17801

$ws.Range("H84").Value = 6568
$ws.Range("I84").Value = 5296.5
$ws.Range("J84").Value = 7839.5
$ws.Range("K84").Value = 52965
$ws.Range("L84").Value = 78395
$ws.Range("M84").Value = -47661
$ws.Range("N84").Value = -89003

$ws.Range("H113").Value = 1407.25
$ws.Range("I113").Value = 654
$ws.Range("K113").Value = 1962
$ws.Range("M113").Value = 208

$ws.Range("H122").Value = 4256.64
$ws.Range("I122").Value = 3510.077
$ws.Range("J122").Value = 5065.4165
$ws.Range("K122").Value = 10530.231
$ws.Range("L122").Value = 15196.2495
$ws.Range("M122").Value = -8080.231
$ws.Range("N122").Value = -20096.2495

$ws.Range("H126").Value = 166666860
$ws.Range("I126").Value = 200000200
$ws.Range("J126").Value = 215
$ws.Range("K126").Value = 600000600
$ws.Range("L126").Value = 645
$ws.Range("M126").Value = -599998130
$ws.Range("N126").Value = -5585

$ws.Range("H132").Value = 25644624
$ws.Range("I132").Value = 66668156
$ws.Range("J132").Value = 4916.125
$ws.Range("K132").Value = 200004468
$ws.Range("L132").Value = 14748.375
$ws.Range("M132").Value = -200001938
$ws.Range("N132").Value = -19808.375

